$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 8997.333000000001
$ws.Range("I10").Value = 8997.333000000001
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 8997.333000000001
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -8704.333000000001
$ws.Range("H32").Value = 3276.75
$ws.Range("I32").Value = 1693.125
$ws.Range("K32").Value = 1693.125
$ws.Range("M32").Value = -1367.125
$ws.Range("H33").Value = 634.2
$ws.Range("I33").Value = 320.2857
$ws.Range("K33").Value = 320.2857
$ws.Range("M33").Value = -91.28570000000002
$ws.Range("H129").Value = 1284.2
$ws.Range("I129").Value = 1105.5
$ws.Range("K129").Value = 3316.5
$ws.Range("M129").Value = 1683.5
$ws.Range("H137").Value = 1039
$ws.Range("I137").Value = 780.1667
$ws.Range("J137").Value = 2074.3333
$ws.Range("K137").Value = 2340.5001
$ws.Range("L137").Value = 6222.999899999999
$ws.Range("M137").Value = 209.4998999999998
$ws.Range("N137").Value = -11322.9999
$ws.Range("H138").Value = 3232.4517
$ws.Range("I138").Value = 1720.0769
$ws.Range("J138").Value = 4324.722
$ws.Range("K138").Value = 5160.2307
$ws.Range("L138").Value = 12974.166
$ws.Range("M138").Value = -20.23070000000007
$ws.Range("N138").Value = -23254.166
$ws.Range("N10").ClearContents()
$ws.Range("N32").ClearContents()
$ws.Range("N33").ClearContents()
$ws.Range("N129").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2867.0442
$ws.Range("I32").Value = 1306.8036
$ws.Range("J32").Value = 10148.167
$ws.Range("K32").Value = 1306.8036
$ws.Range("L32").Value = 10148.167
$ws.Range("M32").Value = -1019.8036
$ws.Range("N32").Value = -10722.167
$ws.Range("H74").Value = 41669484
$ws.Range("I74").Value = 50002180
$ws.Range("K74").Value = 50002180
$ws.Range("M74").Value = -50001306
$ws.Range("H77").Value = 41669484
$ws.Range("I77").Value = 50002180
$ws.Range("K77").Value = 250010900
$ws.Range("M77").Value = -250006532
$ws.Range("H110").Value = 90912810
$ws.Range("I110").Value = 111114664
$ws.Range("J110").Value = 4454
$ws.Range("K110").Value = 111114664
$ws.Range("L110").Value = 4454
$ws.Range("M110").Value = -111112619
$ws.Range("N110").Value = -8544
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 402.5
$ws.Range("I11").Value = 482
$ws.Range("J11").Value = 362.75
$ws.Range("K11").Value = 482
$ws.Range("L11").Value = 362.75
$ws.Range("M11").Value = -342
$ws.Range("N11").Value = -642.75
$ws.Range("H74").Value = 32852
$ws.Range("J74").Value = 32852
$ws.Range("L74").Value = 32852
$ws.Range("N74").Value = -34724
$ws.Range("H77").Value = 32852
$ws.Range("J77").Value = 32852
$ws.Range("L77").Value = 98556
$ws.Range("N77").Value = -107916
$ws.Range("H105").Value = 2045.75
$ws.Range("I105").Value = 1964.0769
$ws.Range("J105").Value = 2399.6667
$ws.Range("K105").Value = 1964.0769
$ws.Range("L105").Value = 2399.6667
$ws.Range("M105").Value = -217.0769
$ws.Range("N105").Value = -5893.6667
$ws.Range("H107").Value = 26316442
$ws.Range("I107").Value = 669.8823
$ws.Range("K107").Value = 669.8823
$ws.Range("M107").Value = 1250.1177
$ws.Range("H134").Value = 1435.5
$ws.Range("I134").Value = 1372.8889
$ws.Range("J134").Value = 1999
$ws.Range("K134").Value = 4118.6667
$ws.Range("L134").Value = 5997
$ws.Range("M134").Value = -1583.6667
$ws.Range("N134").Value = -11067
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("N137").ClearContents()
$ws.Range("M138").ClearContents()
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3135.1594
$ws.Range("I31").Value = 1669.9231
$ws.Range("J31").Value = 3475.3035
$ws.Range("K31").Value = 1669.9231
$ws.Range("L31").Value = 3475.3035
$ws.Range("M31").Value = -1374.9231
$ws.Range("N31").Value = -4065.3035
$ws.Range("H34").Value = 3135.1594
$ws.Range("I34").Value = 1669.9231
$ws.Range("J34").Value = 3475.3035
$ws.Range("K34").Value = 1669.9231
$ws.Range("L34").Value = 3475.3035
$ws.Range("M34").Value = -1467.9231
$ws.Range("N34").Value = -3879.3035
$ws.Range("H58").Value = 1615.6428
$ws.Range("I58").Value = 1374.6364
$ws.Range("K58").Value = 1374.6364
$ws.Range("M58").Value = -1171.6364
$ws.Range("H122").Value = 3031.1538
$ws.Range("I122").Value = 2101.111
$ws.Range("K122").Value = 6303.333
$ws.Range("M122").Value = -3853.333
$ws.Range("H132").Value = 2219.6
$ws.Range("I132").Value = 2096.5757
$ws.Range("K132").Value = 6289.7271
$ws.Range("M132").Value = -3759.7271
$ws.Range("H134").Value = 1056.3334
$ws.Range("I134").Value = 913.96295
$ws.Range("K134").Value = 2741.88885
$ws.Range("M134").Value = -206.8888499999998
$ws.Range("H136").Value = 1615.6428
$ws.Range("I136").Value = 1374.6364
$ws.Range("K136").Value = 4123.9092
$ws.Range("M136").Value = -1573.9092
$ws.Range("N58").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("N134").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2178
$ws.Range("J34").Value = 2436.875
$ws.Range("L34").Value = 7310.625
$ws.Range("N34").Value = -7478.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 7006
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 7006
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 7006
$ws.Range("N13").Value = -7284
$ws.Range("H25").Value = 2424.5
$ws.Range("I25").Value = 2300
$ws.Range("J25").Value = 2466
$ws.Range("K25").Value = 2300
$ws.Range("L25").Value = 2466
$ws.Range("M25").Value = -1771
$ws.Range("N25").Value = -3524
$ws.Range("H28").Value = 25000
$ws.Range("I28").Value = 25000
$ws.Range("K28").Value = 25000
$ws.Range("M28").Value = -24808
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("H132").Value = 3107.8298
$ws.Range("I132").Value = 2932.353
$ws.Range("K132").Value = 8797.059000000001
$ws.Range("M132").Value = -6267.059000000001
$ws.Range("M13").ClearContents()
$ws.Range("M63").ClearContents()
$ws.Range("M66").ClearContents()
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 6000
$ws.Range("I19").Value = 6000
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 6000
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -5830
$ws.Range("H30").Value = 2102.5
$ws.Range("I30").Value = 2005
$ws.Range("J30").Value = 2395
$ws.Range("K30").Value = 2005
$ws.Range("L30").Value = 2395
$ws.Range("M30").Value = -1897
$ws.Range("N30").Value = -2611
$ws.Range("H93").Value = 2032.1666
$ws.Range("I93").Value = 2017.0625
$ws.Range("J93").Value = 2062.375
$ws.Range("K93").Value = 2017.0625
$ws.Range("L93").Value = 2062.375
$ws.Range("M93").Value = -769.0625
$ws.Range("N93").Value = -4558.375
$ws.Range("H132").Value = 6966.1
$ws.Range("I132").Value = 3742.875
$ws.Range("J132").Value = 9114.916999999999
$ws.Range("K132").Value = 11228.625
$ws.Range("L132").Value = 27344.751
$ws.Range("M132").Value = -8698.625
$ws.Range("N132").Value = -32404.751
$ws.Range("N19").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 609.0625
$ws.Range("I107").Value = 470.9565
$ws.Range("J107").Value = 962
$ws.Range("K107").Value = 1412.8695
$ws.Range("L107").Value = 2886
$ws.Range("M107").Value = 507.1305
$ws.Range("N107").Value = -6726
$ws.Range("H132").Value = 6187.4375
$ws.Range("I132").Value = 6100.0835
$ws.Range("K132").Value = 18300.2505
$ws.Range("M132").Value = -15770.2505
$ws.Range("N132").ClearContents()
